# Update SampleRules.xlsx with test data for dynamic scrolling verification
#
# The rule table's third condition column ("CONDITION3" header sits in D,
# but the third generated condition column - "customer.getCreditScore()
# >= 700" - lives in F) is removed entirely. Deleting the whole column
# shifts the two ACTION columns (previously G and H) one place to the
# left (now F and G), and shrinks the used range from A1:H11 to A1:G11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F:F").Delete()
